$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row entirely (row 26), and the "SC 92" row
# (originally row 28, which becomes row 27 once "RM 232" is removed).
# Every following row shifts up by two, matching dimension A1:F35 -> A1:F33.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Update the "missing data" mask: some previously-blank F (error) cells now
# carry a computed value, while some previously-filled ones are now blank.
$ws.Range("F2").Value = 18.03
$ws.Range("F6").ClearContents()
$ws.Range("F12").Value = 17.45
$ws.Range("F14").ClearContents()
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F23").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("F31").Value = 17.18
$ws.Range("F33").Value = 17.53

# A couple of D-column-labelled (E) cells also flip which rows are blanked.
$ws.Range("E26").Value = -5
$ws.Range("E27").ClearContents()
$ws.Range("E30").Value = -5.7
$ws.Range("E32").ClearContents()

# Cells that were already blank and remain blank - re-assert so the saved
# representation is uniform with the other cleared cells above.
$ws.Range("F4").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F22").ClearContents()
